# The workbook's single "Arándano (blue)" price sheet gained a new weekly
# price-report row. It must be inserted as row 87 (pushing every following
# row down by one, old row 87 -> 88, ..., old row 155 -> 156), matching the
# "Fruta / hortaliza, semanal" refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 87, shifting rows 87:155 down
# to 88:156 (row count grows from 155 to 156, matching dimension A1:T156).
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new report's values.
$ws.Range("A87").Value = 9
$ws.Range("B87").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C87").Value = "Metropolitana"
$ws.Range("D87").Value = 44574
$ws.Range("E87").Value = 13
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100101
$ws.Range("H87").Value = "Berries"
$ws.Range("I87").Value = 100101001
$ws.Range("J87").Value = "Arándano (blue)"
$ws.Range("K87").Value = "Sin especificar"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 380
$ws.Range("N87").Value = 4000
$ws.Range("O87").Value = 4000
$ws.Range("P87").Value = 4000
$ws.Range("Q87").Value = "$/bandeja 2 kilos"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 2000
$ws.Range("T87").Value = 2
